$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 92
$ws.Range("H92").Value = 14403920
$ws.Range("I92").Value = 2525850.8
$ws.Range("J92").Value = 66667424
$ws.Range("K92").Value = 2525850.8
$ws.Range("L92").Value = 66667424
$ws.Range("M92").Value = -2524602.8
$ws.Range("N92").Value = -66669920

# Row 103
$ws.Range("H103").Value = 1108
$ws.Range("I103").Value = 772
$ws.Range("J103").Value = 1500
$ws.Range("K103").Value = 2316
$ws.Range("L103").Value = 4500
$ws.Range("M103").Value = -1730
$ws.Range("N103").Value = -5672

# Row 137
$ws.Range("H137").Value = 1684.8096
$ws.Range("I137").Value = 1653.4
$ws.Range("J137").Value = 1763.3334
$ws.Range("K137").Value = 4960.200000000001
$ws.Range("L137").Value = 5290.0002
$ws.Range("M137").Value = -2410.200000000001
$ws.Range("N137").Value = -10390.0002

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4683.136
$ws.Range("I32").Value = 3412.1887
$ws.Range("K32").Value = 3412.1887
$ws.Range("M32").Value = -3125.1887

# Row 45
$ws.Range("H45").Value = 17685.166
$ws.Range("I45").Value = 17685.166
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 17685.166
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -17308.166
$ws.Range("N45").ClearContents()

# Row 61
$ws.Range("H61").Value = 6377.04
$ws.Range("I61").Value = 8814.134
$ws.Range("J61").Value = 2721.4
$ws.Range("K61").Value = 8814.134
$ws.Range("L61").Value = 2721.4
$ws.Range("M61").Value = -8602.134
$ws.Range("N61").Value = -3145.4

# Row 122
$ws.Range("H122").Value = 1606646.9
$ws.Range("I122").Value = 2853711.2
$ws.Range("J122").Value = 3278.4285
$ws.Range("K122").Value = 8561133.600000001
$ws.Range("L122").Value = 9835.2855
$ws.Range("M122").Value = -8558683.600000001
$ws.Range("N122").Value = -14735.2855

# Row 132
$ws.Range("H132").Value = 2733.738
$ws.Range("I132").Value = 1041.5333
$ws.Range("J132").Value = 6964.25
$ws.Range("K132").Value = 3124.5999
$ws.Range("L132").Value = 20892.75
$ws.Range("M132").Value = -594.5999000000002
$ws.Range("N132").Value = -25952.75

# Row 136
$ws.Range("H136").Value = 6377.04
$ws.Range("I136").Value = 8814.134
$ws.Range("J136").Value = 2721.4
$ws.Range("K136").Value = 26442.402
$ws.Range("L136").Value = 8164.200000000001
$ws.Range("M136").Value = -23892.402
$ws.Range("N136").Value = -13264.2

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 76924136
$ws.Range("I99").Value = 83334350
$ws.Range("K99").Value = 83334350
$ws.Range("M99").Value = -83332852

# Row 122
$ws.Range("H122").Value = 57997.5
$ws.Range("J122").Value = 57997.5
$ws.Range("L122").Value = 57997.5
$ws.Range("N122").Value = -67797.5

# Row 134
$ws.Range("H134").Value = 6686.72
$ws.Range("I134").Value = 11097.417
$ws.Range("J134").Value = 2615.3076
$ws.Range("K134").Value = 33292.251
$ws.Range("L134").Value = 7845.9228
$ws.Range("M134").Value = -30757.251
$ws.Range("N134").Value = -12915.9228

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 4808823
$ws.Range("I16").Value = 6411222.5
$ws.Range("K16").Value = 6411222.5
$ws.Range("M16").Value = -6410935.5

# Row 31
$ws.Range("H31").Value = 6484.129
$ws.Range("I31").Value = 1724.4348
$ws.Range("J31").Value = 20168.25
$ws.Range("K31").Value = 1724.4348
$ws.Range("L31").Value = 20168.25
$ws.Range("M31").Value = -1429.4348
$ws.Range("N31").Value = -20758.25

# Row 34
$ws.Range("H34").Value = 6484.129
$ws.Range("I34").Value = 1724.4348
$ws.Range("J34").Value = 20168.25
$ws.Range("K34").Value = 1724.4348
$ws.Range("L34").Value = 20168.25
$ws.Range("M34").Value = -1522.4348
$ws.Range("N34").Value = -20572.25

# Row 58
$ws.Range("H58").Value = 1406.0646
$ws.Range("I58").Value = 1009.7059
$ws.Range("J58").Value = 1887.3572
$ws.Range("K58").Value = 1009.7059
$ws.Range("L58").Value = 1887.3572
$ws.Range("M58").Value = -806.7059
$ws.Range("N58").Value = -2293.3572

# Row 113
$ws.Range("H113").Value = 4808823
$ws.Range("I113").Value = 6411222.5
$ws.Range("K113").Value = 6411222.5
$ws.Range("M113").Value = -6409052.5

# Row 132
$ws.Range("H132").Value = 2627.9312
$ws.Range("I132").Value = 2641.9
$ws.Range("J132").Value = 2596.889
$ws.Range("K132").Value = 7925.700000000001
$ws.Range("L132").Value = 7790.667
$ws.Range("M132").Value = -5395.700000000001
$ws.Range("N132").Value = -12850.667

# Row 136
$ws.Range("H136").Value = 1406.0646
$ws.Range("I136").Value = 1009.7059
$ws.Range("J136").Value = 1887.3572
$ws.Range("K136").Value = 3029.1177
$ws.Range("L136").Value = 5662.071599999999
$ws.Range("M136").Value = -479.1177000000002
$ws.Range("N136").Value = -10762.0716

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 665.5
$ws.Range("J2").Value = 93.59999999999999
$ws.Range("L2").Value = 561.5999999999999
$ws.Range("N2").Value = -787.5999999999999

# Row 5
$ws.Range("H5").Value = 334227.44
$ws.Range("J5").Value = 667900
$ws.Range("L5").Value = 2003700
$ws.Range("N5").Value = -2003924

# Row 107
$ws.Range("H107").Value = 359.93332
$ws.Range("I107").Value = 180
$ws.Range("J107").Value = 449.9
$ws.Range("K107").Value = 540
$ws.Range("L107").Value = 1349.7
$ws.Range("M107").Value = 1380
$ws.Range("N107").Value = -5189.7

# Row 113
$ws.Range("H113").Value = 5000591
$ws.Range("J113").Value = 1250680
$ws.Range("L113").Value = 3752040
$ws.Range("N113").Value = -3756380

# Row 116
$ws.Range("H116").Value = 1951.4
$ws.Range("I116").Value = 585.6667
$ws.Range("K116").Value = 1757.0001
$ws.Range("M116").Value = 1684.9999

# Row 135
$ws.Range("H135").Value = 334227.44
$ws.Range("J135").Value = 667900
$ws.Range("L135").Value = 6011100
$ws.Range("N135").Value = -6016170

$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

# Row 122
$ws.Range("H122").Value = 3088189
$ws.Range("I122").Value = 4987289
$ws.Range("J122").Value = 2151.75
$ws.Range("K122").Value = 14961867
$ws.Range("L122").Value = 6455.25
$ws.Range("M122").Value = -14959417
$ws.Range("N122").Value = -11355.25

$ws = $wb.Worksheets.Item("LTW")
# Row 54
$ws.Range("H54").Value = 13648.333
$ws.Range("J54").Value = 13648.333
$ws.Range("L54").Value = 13648.333
$ws.Range("N54").Value = -14936.333

# Row 136
$ws.Range("H136").Value = 8949.68
$ws.Range("I136").Value = 8610.111000000001
$ws.Range("J136").Value = 9822.857
$ws.Range("K136").Value = 25830.333
$ws.Range("L136").Value = 29468.571
$ws.Range("M136").Value = -23280.333
$ws.Range("N136").Value = -34568.571

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1691.2
$ws.Range("I132").Value = 1018.86664
$ws.Range("J132").Value = 2699.7
$ws.Range("K132").Value = 3056.59992
$ws.Range("L132").Value = 8099.099999999999
$ws.Range("M132").Value = -526.5999199999997
$ws.Range("N132").Value = -13159.1

# Row 136
$ws.Range("H136").Value = 809.53125
$ws.Range("I136").Value = 550.7692
$ws.Range("J136").Value = 1930.8334
$ws.Range("K136").Value = 1652.3076
$ws.Range("L136").Value = 5792.5002
$ws.Range("M136").Value = 897.6924000000001
$ws.Range("N136").Value = -10892.5002
